# Adds the initial "auth" category rows (service-name / category / base-path /
# default-port / source) to the services-details sheet, mirroring the existing
# row layout/format used by the other categories (e.g. the "search" rows which
# use the same fill/border style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (fill colour + border) from an existing styled row so the
# new "auth" rows visually match the rest of the table, then fill in the data.
$ws.Range("A9:E9").Copy()
$ws.Range("A52:E54").PasteSpecial(-4122)

$ws.Cells.Item(52, 1).Value = 54
$ws.Cells.Item(52, 2).Value = "auth"
$ws.Cells.Item(52, 3).Value = "auth-manager"
$ws.Cells.Item(52, 4).Value = 7040
$ws.Cells.Item(52, 5).Value = "/selling/auth/manager/"

$ws.Cells.Item(53, 1).Value = 55
$ws.Cells.Item(53, 2).Value = "auth"
$ws.Cells.Item(53, 3).Value = "auth-users-manager"
$ws.Cells.Item(53, 4).Value = 7041

$ws.Cells.Item(54, 1).Value = 56
$ws.Cells.Item(54, 2).Value = "auth"
$ws.Cells.Item(54, 3).Value = "auth-profiles-manager"
$ws.Cells.Item(54, 4).Value = 7042

$ws.Cells.Item(53, 5).Value = "/selling/auth/users/manager/"
$ws.Cells.Item(54, 5).Value = "/selling/auth/profiles/manager/"

# Match the author's final viewport selection.
$null = $ws.Range("F39").Select()
